# Word COM-interop script: update "Major Project Proposal 2022.docx"
# per the authored diff - rewrites the intro paragraph under "Defining the
# problem" heading and inserts a new paragraph under "Objectives and design
# specifications".

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Step 1: Replace the body paragraph under the "Defining the problem" heading
# (paragraph 3) with the new write-up text, preserving the exact run /
# proofErr-span structure of the authored edit (e.g. spell-check spans around
# "lootbox", "Skinport", "Skinbaron").
# ---------------------------------------------------------------------------
$introXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Within the game Counter-Strike: Global Offensive (CS:GO) there exists a “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>lootbox</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">” system where players pay real money for a case and a key to roll for a random item contained within, with the potential of receiving </w:t></w:r><w:r><w:t>a rare item worth lots of money on</w:t></w:r><w:r><w:t xml:space="preserve"> the</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>in-game</w:t></w:r><w:r><w:t xml:space="preserve"> marketplace </w:t></w:r><w:r><w:t>called the “</w:t></w:r><w:r><w:t>Steam Market</w:t></w:r><w:r><w:t>”</w:t></w:r><w:r><w:t xml:space="preserve"> or external sites such as “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Skinport</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>” or “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Skinbaron</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">” </w:t></w:r><w:r><w:t xml:space="preserve">which allow for users to cash out their items for real money by selling them to other players or the marketplace itself. </w:t></w:r><w:r><w:t>This system presents the curiosity, if one were to have a large sum of money, how profitable would gambling it on CS:GO be? Which choice of case is most profitable</w:t></w:r><w:r><w:t xml:space="preserve"> long term</w:t></w:r><w:r><w:t>?</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">However, there is no capability within the game to simulate such an experience without the risk </w:t></w:r><w:r><w:t>of spending copious amounts of your own money.</w:t></w:r><w:r><w:t xml:space="preserve"> My solution </w:t></w:r><w:r><w:t>was</w:t></w:r><w:r><w:t xml:space="preserve"> to create a price accurate simulation of the </w:t></w:r><w:r><w:t>case opening system within CS:GO</w:t></w:r><w:r><w:t>, allowing users to purchase cases and sell items much alike to the game</w:t></w:r><w:r><w:t>, without the aspect of spending actual money.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$problemPara = $d.Paragraphs.Item(3)
$problemRange = $problemPara.Range
$problemRange.End = $problemRange.End - 1
$problemRange.Text = ""
$problemRange.InsertXML($introXml)

# ---------------------------------------------------------------------------
# Step 2: Insert a brand-new paragraph directly after the "Objectives and
# design specifications" heading (paragraph 4), ahead of the pre-existing
# "main feature objective" paragraph (originally paragraph 5).
# ---------------------------------------------------------------------------
$objectivesXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">The idea comprises a website that simulates the case opening and item system in the game Counter-Strike: Global Offensive (CS:GO). The simulation would provide the user with a </w:t></w:r><w:r><w:t xml:space="preserve">choice for their </w:t></w:r><w:r><w:t xml:space="preserve">starting amount of money and allow them to purchase cases much like within the game. A case has a range of potential items contained within with varying value attached to each item, when a case is rolled a single item from within its contents is received. The odds of receiving each item in the case will be provided to the user, </w:t></w:r><w:r><w:t>higher rarity items such as “Souvenir Weapons”, “Knives” or “Gloves” will have much higher value than other items</w:t></w:r><w:r><w:t xml:space="preserve">. The simulation will allow the user to sell the items they receive to increase their balance and open more cases. Each case will have a different price depending on the rarity of the items within and the age of the case itself, as well as an accurate average market price for each item which will be scraped from an active source each time the site is loaded. </w:t></w:r><w:r><w:t>The system will track the users wins and losses on case openings and allow them to look at their statistics once they either choose to end the simulation or run themselves completely out of money.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$followingPara = $d.Paragraphs.Item(5)
$followingPara.Range.InsertParagraphBefore()

$newPara = $d.Paragraphs.Item(5)
$newRange = $newPara.Range
$newRange.End = $newRange.End - 1
$newRange.InsertXML($objectivesXml)
